$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "aGKBW286"
$ws.Range("B2").Value = 23111603
$ws.Range("C2").Value = "nhfldgk95"
$ws.Range("D2").Value = "wr3V9!`$M"
$ws.Range("F2").Value = "ASiXJWbN"
$ws.Range("G2").Value = "WKRh"
